$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 3432
$ws1.Range("F3").Value = 744
$ws1.Range("F5").Value = 6980
$ws1.Range("F6").Value = 2475
$ws1.Range("F11").Value = 78
$ws1.Range("F13").Value = 174
$ws1.Range("F14").Value = 577

# Sheet "演出" (Performance)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 23

# Sheet "全部类型" (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 3432
$ws4.Range("F3").Value = 23
$ws4.Range("F4").Value = 744
$ws4.Range("F6").Value = 6980
$ws4.Range("F7").Value = 2475
$ws4.Range("F12").Value = 78
$ws4.Range("F14").Value = 174
$ws4.Range("F15").Value = 577
